# Update the division-problem table: replace each equation's old
# "A÷B=" text with the new value, cell by cell (row, column) so that
# duplicate equation strings elsewhere in the table are left untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $rng = $cell.Range
    # A table-cell Range includes the trailing end-of-cell marker, so
    # trim the last character before replacing the visible text.
    $target = $d.Range($rng.Start, $rng.End - 1)
    $target.Text = $newText
}

$updates = @(
    @{ Row = 1;  Col = 1; Text = "58÷6=" },
    @{ Row = 1;  Col = 2; Text = "14÷4=" },
    @{ Row = 1;  Col = 3; Text = "89÷8=" },
    @{ Row = 1;  Col = 4; Text = "83÷8=" },
    @{ Row = 1;  Col = 5; Text = "66÷5=" },

    @{ Row = 5;  Col = 1; Text = "91÷6=" },
    @{ Row = 5;  Col = 2; Text = "25÷6=" },
    @{ Row = 5;  Col = 3; Text = "52÷4=" },
    @{ Row = 5;  Col = 4; Text = "73÷6=" },
    @{ Row = 5;  Col = 5; Text = "15÷5=" },

    @{ Row = 9;  Col = 1; Text = "22÷4=" },
    @{ Row = 9;  Col = 2; Text = "78÷8=" },
    @{ Row = 9;  Col = 3; Text = "69÷3=" },
    @{ Row = 9;  Col = 4; Text = "98÷5=" },
    @{ Row = 9;  Col = 5; Text = "61÷2=" },

    @{ Row = 13; Col = 1; Text = "72÷5=" },
    @{ Row = 13; Col = 2; Text = "80÷8=" },
    @{ Row = 13; Col = 3; Text = "27÷3=" },
    @{ Row = 13; Col = 4; Text = "37÷7=" },
    @{ Row = 13; Col = 5; Text = "48÷8=" },

    @{ Row = 17; Col = 1; Text = "39÷9=" },
    @{ Row = 17; Col = 2; Text = "53÷5=" },
    @{ Row = 17; Col = 3; Text = "91÷2=" },
    @{ Row = 17; Col = 4; Text = "22÷2=" },
    @{ Row = 17; Col = 5; Text = "62÷3=" }
)

foreach ($u in $updates) {
    Set-CellText $t $u.Row $u.Col $u.Text
}

Write-Host "Done updating" $updates.Count "cells"
